$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 505.08334
$ws.Range("I2").Value = 298.7143
$ws.Range("J2").Value = 794
$ws.Range("K2").Value = 298.7143
$ws.Range("L2").Value = 794
$ws.Range("M2").Value = -185.7143
$ws.Range("N2").Value = -1020
$ws.Range("H29").Value = 3430
$ws.Range("I29").Value = 50
$ws.Range("K29").Value = 150
$ws.Range("M29").Value = 131
$ws.Range("H42").Value = 1749.6923
$ws.Range("I42").Value = 1154
$ws.Range("J42").Value = 3735.3333
$ws.Range("K42").Value = 3462
$ws.Range("L42").Value = 11205.9999
$ws.Range("M42").Value = -3232
$ws.Range("N42").Value = -11665.9999
$ws.Range("H70").Value = 3069.0938
$ws.Range("J70").Value = 6416.4165
$ws.Range("L70").Value = 19249.2495
$ws.Range("N70").Value = -19789.2495
$ws.Range("H73").Value = 3069.0938
$ws.Range("J73").Value = 6416.4165
$ws.Range("L73").Value = 19249.2495
$ws.Range("N73").Value = -21121.2495
$ws.Range("H125").Value = 899587.8
$ws.Range("I125").Value = 3392546.5
$ws.Range("J125").Value = 9245.357
$ws.Range("K125").Value = 30532918.5
$ws.Range("L125").Value = 83208.213
$ws.Range("M125").Value = -30530458.5
$ws.Range("N125").Value = -88128.213

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3613.5
$ws.Range("I32").Value = 2357.024
$ws.Range("K32").Value = 2357.024
$ws.Range("M32").Value = -2070.024
$ws.Range("H74").Value = 22731436
$ws.Range("I74").Value = 25003954
$ws.Range("K74").Value = 25003954
$ws.Range("M74").Value = -25003080
$ws.Range("H75").Value = 150000
$ws.Range("J75").Value = 150000
$ws.Range("L75").Value = 150000
$ws.Range("N75").Value = -151748
$ws.Range("H77").Value = 22731436
$ws.Range("I77").Value = 25003954
$ws.Range("K77").Value = 125019770
$ws.Range("M77").Value = -125015402
$ws.Range("H78").Value = 150000
$ws.Range("J78").Value = 150000
$ws.Range("L78").Value = 450000
$ws.Range("N78").Value = -458736
$ws.Range("H132").Value = 6255135.5
$ws.Range("I132").Value = 7147369
$ws.Range("J132").Value = 9500
$ws.Range("K132").Value = 21442107
$ws.Range("L132").Value = 28500
$ws.Range("M132").Value = -21439577
$ws.Range("N132").Value = -33560

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 89208.914
$ws.Range("I107").Value = 5877.25
$ws.Range("J107").Value = 255872.25
$ws.Range("K107").Value = 5877.25
$ws.Range("L107").Value = 255872.25
$ws.Range("M107").Value = -3957.25
$ws.Range("N107").Value = -259712.25
$ws.Range("H134").Value = 18521858
$ws.Range("I134").Value = 20003192
$ws.Range("K134").Value = 60009576
$ws.Range("M134").Value = -60007041

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 732986.25
$ws.Range("I16").Value = 998072.5600000001
$ws.Range("J16").Value = 3999
$ws.Range("K16").Value = 998072.5600000001
$ws.Range("L16").Value = 3999
$ws.Range("M16").Value = -997785.5600000001
$ws.Range("N16").Value = -4573
$ws.Range("H31").Value = 4209
$ws.Range("I31").Value = 4349.875
$ws.Range("J31").Value = 3833.3333
$ws.Range("K31").Value = 4349.875
$ws.Range("L31").Value = 3833.3333
$ws.Range("M31").Value = -4054.875
$ws.Range("N31").Value = -4423.3333
$ws.Range("H34").Value = 4209
$ws.Range("I34").Value = 4349.875
$ws.Range("J34").Value = 3833.3333
$ws.Range("K34").Value = 4349.875
$ws.Range("L34").Value = 3833.3333
$ws.Range("M34").Value = -4147.875
$ws.Range("N34").Value = -4237.3333
$ws.Range("H51").Value = 45000
$ws.Range("J51").Value = 45000
$ws.Range("L51").Value = 45000
$ws.Range("N51").Value = -46472
$ws.Range("H61").Value = 45000
$ws.Range("J61").Value = 45000
$ws.Range("L61").Value = 45000
$ws.Range("N61").Value = -45696
$ws.Range("H86").Value = 11584.546
$ws.Range("I86").Value = 9060.362999999999
$ws.Range("J86").Value = 14108.728
$ws.Range("K86").Value = 9060.362999999999
$ws.Range("L86").Value = 14108.728
$ws.Range("M86").Value = -7937.362999999999
$ws.Range("N86").Value = -16354.728
$ws.Range("H89").Value = 11584.546
$ws.Range("I89").Value = 9060.362999999999
$ws.Range("J89").Value = 14108.728
$ws.Range("K89").Value = 45301.815
$ws.Range("L89").Value = 70543.64
$ws.Range("M89").Value = -39685.815
$ws.Range("N89").Value = -81775.64
$ws.Range("H113").Value = 732986.25
$ws.Range("I113").Value = 998072.5600000001
$ws.Range("J113").Value = 3999
$ws.Range("K113").Value = 998072.5600000001
$ws.Range("L113").Value = 3999
$ws.Range("M113").Value = -995902.5600000001
$ws.Range("N113").Value = -8339
$ws.Range("H122").Value = 3089.375
$ws.Range("I122").Value = 1940
$ws.Range("K122").Value = 5820
$ws.Range("M122").Value = -3370
$ws.Range("H134").Value = 50205132
$ws.Range("I134").Value = 62754470
$ws.Range("K134").Value = 188263410
$ws.Range("M134").Value = -188260875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 6249.231
$ws.Range("I3").Value = 4658.273
$ws.Range("J3").Value = 14999.5
$ws.Range("K3").Value = 13974.819
$ws.Range("L3").Value = 44998.5
$ws.Range("M3").Value = -13862.819
$ws.Range("N3").Value = -45222.5
$ws.Range("H7").Value = 2007798.4
$ws.Range("J7").Value = 9501
$ws.Range("L7").Value = 28503
$ws.Range("N7").Value = -28727

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4450
$ws.Range("I70").Value = 4175
$ws.Range("K70").Value = 4175
$ws.Range("M70").Value = -3905
$ws.Range("H73").Value = 4450
$ws.Range("I73").Value = 4175
$ws.Range("K73").Value = 4175
$ws.Range("M73").Value = -3239
$ws.Range("H113").Value = 62603.293
$ws.Range("I113").Value = 70317.07000000001
$ws.Range("J113").Value = 4750
$ws.Range("K113").Value = 70317.07000000001
$ws.Range("L113").Value = 4750
$ws.Range("M113").Value = -68147.07000000001
$ws.Range("N113").Value = -9090
$ws.Range("H122").Value = 70279.05499999999
$ws.Range("I122").Value = 82535
$ws.Range("J122").Value = 8999.333000000001
$ws.Range("K122").Value = 247605
$ws.Range("L122").Value = 26997.999
$ws.Range("M122").Value = -245155
$ws.Range("N122").Value = -31897.999
$ws.Range("H127").Value = 65000
$ws.Range("J127").Value = 65000
$ws.Range("L127").Value = 65000
$ws.Range("N127").Value = -74920
$ws.Range("H132").Value = 9622311
$ws.Range("I132").Value = 11370013
$ws.Range("K132").Value = 34110039
$ws.Range("M132").Value = -34107509

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1321.2222
$ws.Range("I46").Value = 1356.1428
$ws.Range("K46").Value = 1356.1428
$ws.Range("M46").Value = -1168.1428
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("H55").Value = 412.07693
$ws.Range("I55").Value = 386.8
$ws.Range("J55").Value = 496.33334
$ws.Range("K55").Value = 386.8
$ws.Range("L55").Value = 496.33334
$ws.Range("M55").Value = -213.8
$ws.Range("N55").Value = -842.33334
$ws.Range("H68").Value = 5265355.5
$ws.Range("I68").Value = 8773761
$ws.Range("J68").Value = 2747.5
$ws.Range("K68").Value = 8773761
$ws.Range("L68").Value = 2747.5
$ws.Range("M68").Value = -8773012
$ws.Range("N68").Value = -4245.5
$ws.Range("H71").Value = 5265355.5
$ws.Range("I71").Value = 8773761
$ws.Range("J71").Value = 2747.5
$ws.Range("K71").Value = 43868805
$ws.Range("L71").Value = 13737.5
$ws.Range("M71").Value = -43865061
$ws.Range("N71").Value = -21225.5
$ws.Range("H82").Value = 964.2593000000001
$ws.Range("I82").Value = 1062.65
$ws.Range("J82").Value = 683.1429000000001
$ws.Range("K82").Value = 1062.65
$ws.Range("L82").Value = 683.1429000000001
$ws.Range("M82").Value = -701.6500000000001
$ws.Range("N82").Value = -1405.1429
$ws.Range("H85").Value = 964.2593000000001
$ws.Range("I85").Value = 1062.65
$ws.Range("J85").Value = 683.1429000000001
$ws.Range("K85").Value = 1062.65
$ws.Range("L85").Value = 683.1429000000001
$ws.Range("M85").Value = 185.3499999999999
$ws.Range("N85").Value = -3179.1429
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H122").Value = 6131.8887
$ws.Range("I122").Value = 6131.8887
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 18395.6661
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -15945.6661
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 22736360
$ws.Range("I132").Value = 22736360
$ws.Range("K132").Value = 68209080
$ws.Range("M132").Value = -68206550
$ws.Range("H136").Value = 1082.8334
$ws.Range("I136").Value = 904.7
$ws.Range("J136").Value = 1973.5
$ws.Range("K136").Value = 2714.1
$ws.Range("L136").Value = 5920.5
$ws.Range("M136").Value = -164.1000000000004
$ws.Range("N136").Value = -11020.5
$ws.Range("H137").Value = 85000
$ws.Range("J137").Value = 85000
$ws.Range("L137").Value = 85000
$ws.Range("N137").Value = -95200

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 843.125
$ws.Range("I100").Value = 763.4091
$ws.Range("J100").Value = 1720
$ws.Range("K100").Value = 1526.8182
$ws.Range("L100").Value = 3440
$ws.Range("M100").Value = -985.8181999999999
$ws.Range("N100").Value = -4522
$ws.Range("H132").Value = 13164822
$ws.Range("I132").Value = 17860166
$ws.Range("K132").Value = 53580498
$ws.Range("M132").Value = -53577968
